$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 3196
$ws1.Range("F11").Value = 8131
$ws1.Range("F15").Value = 23
$ws1.Range("F16").Value = 290
$ws1.Range("F20").Value = 313
$ws1.Range("F21").Value = 10151
$ws1.Range("F29").Value = 141
$ws1.Range("F30").Value = 83
$ws1.Range("F31").Value = 125
$ws1.Range("F32").Value = 42
$ws1.Range("F33").Value = 2061
$ws1.Range("F36").Value = 2102
$ws1.Range("F37").Value = 4032
$ws1.Range("F38").Value = 256
$ws1.Range("F39").Value = 64
$ws1.Range("F40").Value = 2263
$ws1.Range("F41").Value = 1216
$ws1.Range("F42").Value = 145
$ws1.Range("F43").Value = 301
$ws1.Range("F44").Value = 225
$ws1.Range("F45").Value = 27
$ws1.Range("F46").Value = 90
$ws1.Range("F47").Value = 81
$ws1.Range("F48").Value = 78

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 173
$ws2.Range("F6").Value = 39
$ws2.Range("F14").Value = 6
$ws2.Range("F18").Value = 24
$ws2.Range("F19").Value = 37

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 173
$ws4.Range("F8").Value = 3196
$ws4.Range("F13").Value = 8131
$ws4.Range("F17").Value = 23
$ws4.Range("F18").Value = 290
$ws4.Range("F21").Value = 313
$ws4.Range("F22").Value = 10151
$ws4.Range("F29").Value = 141
$ws4.Range("F30").Value = 83
$ws4.Range("F31").Value = 125
$ws4.Range("F32").Value = 42
$ws4.Range("F33").Value = 2061
$ws4.Range("F35").Value = 2102
$ws4.Range("F36").Value = 4032
$ws4.Range("F37").Value = 256
$ws4.Range("F38").Value = 64
$ws4.Range("F39").Value = 2263
$ws4.Range("F40").Value = 24
$ws4.Range("F41").Value = 1216
$ws4.Range("F42").Value = 145
$ws4.Range("F43").Value = 301
$ws4.Range("F44").Value = 225
$ws4.Range("F45").Value = 27
$ws4.Range("F46").Value = 90
$ws4.Range("F47").Value = 81
$ws4.Range("F48").Value = 78
